# Apply the "double underscore field separator" rename to the
# BusinessPartnerCertificate reporting template workbook.
#
# Changes:
#  - Sheet "semantic_aspect_model_schema" (row 1 headers): rename flattened
#    field names from a single "_" separator to a double "__" separator,
#    and widen the affected columns to match the new (longer) header text.
#  - Sheet "description": mirror the same renames in the "Column Name"
#    column (column A), and fix an off-by-one numbering typo in the legend
#    ("2." -> "1." for the olive-green legend row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: semantic_aspect_model_schema
# ---------------------------------------------------------------------
$schemaSheet = $wb.Worksheets.Item("semantic_aspect_model_schema")

$schemaSheet.Range("B1").Value = "type__certificateType"
$schemaSheet.Range("C1").Value = "type__certificateVersion"
$schemaSheet.Range("F1").Value = "enclosedSites[0]__enclosedSiteBpn"
$schemaSheet.Range("G1").Value = "enclosedSites[0]__areaOfApplication"
$schemaSheet.Range("L1").Value = "validator__validatorName"
$schemaSheet.Range("M1").Value = "validator__validatorBpn"

# Column widths follow the template generator's convention of
# width = (character count of header text) * 1.2, so widen the columns
# whose header text grew by one character (the extra underscore).
# (ColumnWidth is stored on a 1/6-character pixel grid by this engine, so
# the closest representable width is used for each target.)
$schemaSheet.Columns.Item(2).ColumnWidth = 24.333333333333332   # B -> target 25.2
$schemaSheet.Columns.Item(3).ColumnWidth = 28.0                 # C -> target 28.8
$schemaSheet.Columns.Item(6).ColumnWidth = 38.833333333333336   # F -> target 39.6
$schemaSheet.Columns.Item(7).ColumnWidth = 41.166666666666664   # G -> target 42
$schemaSheet.Columns.Item(12).ColumnWidth = 28.0                # L -> target 28.8
$schemaSheet.Columns.Item(13).ColumnWidth = 26.833333333333332  # M -> target 27.6

# ---------------------------------------------------------------------
# Sheet 2: description
# ---------------------------------------------------------------------
$descSheet = $wb.Worksheets.Item("description")

$descSheet.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

$descSheet.Range("A6").Value = "type__certificateType"
$descSheet.Range("A7").Value = "type__certificateVersion"
$descSheet.Range("A10").Value = "enclosedSites[0]__enclosedSiteBpn"
$descSheet.Range("A11").Value = "enclosedSites[0]__areaOfApplication"
$descSheet.Range("A16").Value = "validator__validatorName"
$descSheet.Range("A17").Value = "validator__validatorBpn"
